$wb = $excel.ActiveWorkbook

# --- Table-6.1 sheet: fill in the answers grid ---
$ws6 = $wb.Worksheets.Item("Table-6.1")

# Header row (C4:E4) - rename generic "C1/C2/C3" headers to country-specific labels
$ws6.Range("C4").Value = "C1-GBR"
$ws6.Range("D4").Value = "C2-India"
$ws6.Range("E4").Value = "C3-USA"

# Row 5: Total number of Investments (count)
$ws6.Range("C5").Value = 2303
$ws6.Range("D5").Value = 992
$ws6.Range("E5").Value = 38372

# Row 6: Total amount of investment (USD)
$ws6.Range("C6").Value = 75464027
$ws6.Range("D6").Value = 148288951
$ws6.Range("E6").Value = 100160893

# Row 7: Top Sector name (no. of investment-wise)
$ws6.Range("C7").Value = "Others"
$ws6.Range("D7").Value = "Others"
$ws6.Range("E7").Value = "Others"

# Row 8: Second Sector name (no. of investment-wise)
$ws6.Range("C8").Value = "Social..Finance..Analytics..Advertising"
$ws6.Range("D8").Value = "Social..Finance..Analytics..Advertising"
$ws6.Range("E8").Value = "Cleantech...Semiconductors"

# Row 9: Third Sector name (no. of investment-wise)
$ws6.Range("C9").Value = "Cleantech...Semiconductors"
$ws6.Range("D9").Value = "News..Search.and.Messaging"
$ws6.Range("E9").Value = "Social..Finance..Analytics..Advertising"

# Row 10: Number of investments in top sector (3)
$ws6.Range("C10").Value = 580
$ws6.Range("D10").Value = 332
$ws6.Range("E10").Value = 8768

# Row 11: Number of investments in second sector (4)
$ws6.Range("C11").Value = 481
$ws6.Range("D11").Value = 193
$ws6.Range("E11").Value = 8270

# Row 12: Number of investments in third sector (5)
$ws6.Range("C12").Value = 466
$ws6.Range("D12").Value = 154
$ws6.Range("E12").Value = 7825

# Row 13: For point 3 (top sector count-wise), which company received the highest investment?
$ws6.Range("C13").Value = "OneWeb"
$ws6.Range("D13").Value = "Flipkart"
$ws6.Range("E13").Value = "SoFi"

# Row 14: For point 4 (second best sector count-wise), which company received the highest investment?
$ws6.Range("C14").Value = "Liquid Telecom"
$ws6.Range("D14").Value = "Flipkart"
$ws6.Range("E14").Value = "CoreExpress"

# Autofit columns C:E now that they contain real data
$ws6.Columns("C:E").AutoFit() | Out-Null

# Update the active selection / active sheet state to match the saved workbook
$ws6.Range("D20").Select()

$ws1 = $wb.Worksheets.Item("Table -1.1")
$ws1.Range("D7").Select()

$ws6.Activate()
